$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item #26 "Break Away Headers - Straight (40 pin)" (row 28): quantity 1 -> 2
$ws.Range("D28").Value = 2

# New item #29 "Jumper Wires - Connected 6" (F/F, 20 pack)" (row 31)
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 'Jumper Wires - Connected 6" (F/F, 20 pack)'
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 1.95
$ws.Range("F31").Formula = "=E31*D31"
$ws.Range("G31").Value = "https://www.sparkfun.com/products/12796"
$ws.Hyperlinks.Add($ws.Range("G31"), "https://www.sparkfun.com/products/12796") | Out-Null
$ws.Range("G31").Style = $ws.Range("G30").Style

# Keep the running total formula covering the new last data row
$ws.Range("F32").Formula = "=SUM(F3:F31)"

# Restore the saved scroll position / active selection
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D32").Select()
